# Applies the "data refresh" edit described by the commit diff:
#  - Vancouver rows (2-13): collapse the per-row "Name" values down to the
#    two repeating circuit labels (5L51,5L52 / 2L112), and make the odd
#    ("2L112") rows carry the Richmond substation's old Lat/Lng pair.
#  - Richmond rows (14-19): company renamed Hydro Richmond -> Hydro One,
#    name collapsed to "Fort Frances East, Fort Frances West", PP cycles
#    1-2-3 instead of 1-2-3-4, Year now increments every row (2014..2019),
#    and Lat/Lng point at the new Fort Frances coordinates.
#  - The old continuation rows for Richmond (20-25) are deleted outright,
#    which pulls the already-blank trailing rows up into their place and
#    shrinks the used range from J31 to J25.
#  - Column D is widened to fit the longer "Fort Frances East, ..." text.
#  - A new bold 8pt Arial "label" style (left/center, wrap) is stamped on
#    the now-empty F25 cell, mirroring the new cellXfs entry added upstream.
#  - The view is scrolled/selected to match the saved sheetView state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column width: D gets its own (wider) width, split off from D:E ---
$ws.Columns.Item(4).ColumnWidth = 34.6640625

# --- Vancouver block (rows 2-13): repeating Name + paired Lat/Lng ---
$vanRows = 2..13
foreach ($r in $vanRows) {
    if ((($r - 2) % 2) -eq 0) {
        $ws.Cells.Item($r, 4).Value = "5L51,5L52"
        $ws.Cells.Item($r, 7).Value = 49.015546999999998
        $ws.Cells.Item($r, 8).Value = -122.759846
    } else {
        $ws.Cells.Item($r, 4).Value = "2L112"
        $ws.Cells.Item($r, 7).Value = 49.035527000000002
        $ws.Cells.Item($r, 8).Value = -117.901895
    }
}

# --- Richmond block (rows 14-19): new company/name/year/pp/lat/lng ---
$years = @(2014, 2015, 2016, 2017, 2018, 2019)
$pps = @(1, 2, 3, 1, 2, 3)
for ($i = 0; $i -lt 6; $i++) {
    $r = 14 + $i
    $ws.Cells.Item($r, 2).Value = $years[$i]
    $ws.Cells.Item($r, 3).Value = "Hydro One"
    $ws.Cells.Item($r, 4).Value = "Fort Frances East, Fort Frances West"
    $ws.Cells.Item($r, 5).Value = $pps[$i]
    $ws.Cells.Item($r, 7).Value = 48.603104999999999
    $ws.Cells.Item($r, 8).Value = -93.429597999999999
}

# --- drop the old continuation rows (what used to be years 2017-2019
#     duplicated); this shifts the already-blank rows 26-31 up to 20-25
#     and shrinks the sheet dimension to J25 ---
$ws.Rows("20:25").Delete()

# --- new label style applied to F25 (empty cell, bold 8pt Arial black,
#     left/center aligned, wrapped) ---
$f25 = $ws.Range("F25")
$f25.Font.Name = "Arial"
$f25.Font.Size = 8
$f25.Font.Bold = $true
$f25.Font.Color = 0
$f25.HorizontalAlignment = -4131
$f25.VerticalAlignment = -4108
$f25.WrapText = $true

# --- restore the saved view state (scrolled down one row, E24 selected) ---
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("E24").Select()
